# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1132
$ws1.Range("F7").Value = 239
$ws1.Range("F9").Value = 1021
$ws1.Range("F14").Value = 12841
$ws1.Range("F16").Value = 5275
$ws1.Range("F17").Value = 5534

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1132
$ws4.Range("F7").Value = 239
$ws4.Range("F9").Value = 1021
$ws4.Range("F14").Value = 12841
$ws4.Range("F18").Value = 5275
$ws4.Range("F19").Value = 5534
